$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 23.78418466666666
$ws.Range("H2").Value = 71.352554
$ws.Range("I2").Value = 0.06460357633592957
$ws.Range("J2").Value = 0.06460357633592959
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 3998.355177577447
$ws.Range("R2").Value = 35985.19659819703
$ws.Range("S2").Value = 0.01927888463071305
$ws.Range("T2").Value = 0.01927888463071306
$ws.Range("G3").Value = 23.78418466666666
$ws.Range("H3").Value = 71.352554
$ws.Range("I3").Value = 0.06460357633592957
$ws.Range("J3").Value = 0.06460357633592959
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 3876.970410914186
$ws.Range("R3").Value = 34892.73369822768
$ws.Range("S3").Value = 0.01869360323161411
$ws.Range("T3").Value = 0.01869360323161412
$ws.Range("G4").Value = 23.78418466666666
$ws.Range("H4").Value = 71.352554
$ws.Range("I4").Value = 0.06460357633592957
$ws.Range("J4").Value = 0.06460357633592959
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 3948.020770991873
$ws.Range("R4").Value = 35532.18693892686
$ws.Range("S4").Value = 0.01903618702771856
$ws.Range("T4").Value = 0.01903618702771856
$ws.Range("G5").Value = 23.78418466666666
$ws.Range("H5").Value = 71.352554
$ws.Range("I5").Value = 0.06460357633592957
$ws.Range("J5").Value = 0.06460357633592959
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 1575.148879254274
$ws.Range("R5").Value = 14176.33991328846
$ws.Range("S5").Value = 0.007594901445883845
$ws.Range("T5").Value = 0.007594901445883847
$ws.Range("I6").Value = 0.3773880863345054
$ws.Range("J6").Value = 0.3773880863345054
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 23356.78138165882
$ws.Range("R6").Value = 210211.0324349294
$ws.Range("S6").Value = 0.1126194831632276
$ws.Range("T6").Value = 0.1126194831632276
$ws.Range("I7").Value = 0.3773880863345054
$ws.Range("J7").Value = 0.3773880863345054
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.1092005048387987
$ws.Range("T7").Value = 0.1092005048387987
$ws.Range("I8").Value = 0.3773880863345054
$ws.Range("J8").Value = 0.3773880863345054
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 23062.74804085214
$ws.Range("R8").Value = 207564.7323676692
$ws.Range("S8").Value = 0.1112017414661456
$ws.Range("T8").Value = 0.1112017414661456
$ws.Range("I9").Value = 0.3773880863345054
$ws.Range("J9").Value = 0.3773880863345054
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 9201.385665442012
$ws.Range("R9").Value = 82812.47098897809
$ws.Range("S9").Value = 0.04436635686633354
$ws.Range("T9").Value = 0.04436635686633354
$ws.Range("G10").Value = 150.629115
$ws.Range("H10").Value = 451.887345
$ws.Range("I10").Value = 0.4091449703110563
$ws.Range("J10").Value = 0.4091449703110563
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 25322.23451402281
$ws.Range("R10").Value = 227900.1106262053
$ws.Range("S10").Value = 0.1220963161365496
$ws.Range("T10").Value = 0.1220963161365496
$ws.Range("G11").Value = 150.629115
$ws.Range("H11").Value = 451.887345
$ws.Range("I11").Value = 0.4091449703110563
$ws.Range("J11").Value = 0.4091449703110563
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 24553.48501795143
$ws.Range("R11").Value = 220981.3651615629
$ws.Range("S11").Value = 0.1183896337167906
$ws.Range("T11").Value = 0.1183896337167906
$ws.Range("G12").Value = 150.629115
$ws.Range("H12").Value = 451.887345
$ws.Range("I12").Value = 0.4091449703110563
$ws.Range("J12").Value = 0.4091449703110563
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 25003.45851962595
$ws.Range("R12").Value = 225031.1266766335
$ws.Range("S12").Value = 0.1205592726909142
$ws.Range("T12").Value = 0.1205592726909142
$ws.Range("G13").Value = 150.629115
$ws.Range("H13").Value = 451.887345
$ws.Range("I13").Value = 0.4091449703110563
$ws.Range("J13").Value = 0.4091449703110563
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 9975.674381970115
$ws.Range("R13").Value = 89781.06943773103
$ws.Range("S13").Value = 0.0480997477668019
$ws.Range("T13").Value = 0.0480997477668019
$ws.Range("G14").Value = 54.80491966666667
$ws.Range("H14").Value = 164.414759
$ws.Range("I14").Value = 0.1488633670185088
$ws.Range("J14").Value = 0.1488633670185088
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 9213.245582180538
$ws.Range("R14").Value = 82919.21023962484
$ws.Range("S14").Value = 0.04442354187276172
$ws.Range("T14").Value = 0.04442354187276173
$ws.Range("G15").Value = 54.80491966666667
$ws.Range("H15").Value = 164.414759
$ws.Range("I15").Value = 0.1488633670185088
$ws.Range("J15").Value = 0.1488633670185088
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 8933.543650877402
$ws.Range("R15").Value = 80401.89285789662
$ws.Range("S15").Value = 0.04307490198833606
$ws.Range("T15").Value = 0.04307490198833606
$ws.Range("G16").Value = 54.80491966666667
$ws.Range("H16").Value = 164.414759
$ws.Range("I16").Value = 0.1488633670185088
$ws.Range("J16").Value = 0.1488633670185088
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 9097.262076836423
$ws.Range("R16").Value = 81875.35869152781
$ws.Range("S16").Value = 0.04386430375626462
$ws.Range("T16").Value = 0.04386430375626462
$ws.Range("G17").Value = 54.80491966666667
$ws.Range("H17").Value = 164.414759
$ws.Range("I17").Value = 0.1488633670185088
$ws.Range("J17").Value = 0.1488633670185088
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 3629.550854924009
$ws.Range("R17").Value = 32665.95769431608
$ws.Range("S17").Value = 0.01750061940114638
$ws.Range("T17").Value = 0.01750061940114638
